# Add the new "pelada" results rows (440-461) to the Jogadores sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @('Guinha', 2, 2, 2, 0, 1, 0, 0, 0, 0),
    @('Dogão', 2, 2, 2, 2, 1, 0, 0, 0, 0),
    @('Romario', 2, 2, 2, 1, 1, 0, 0, 0, 0),
    @('Fernando', 2, 2, 2, 2, 1, 0, 0, 1, 0),
    @('Euler', 2, 2, 2, 0, 1, 0, 0, 0, 0),
    @('Peixe', 2, 2, 2, 1, 1, 0, 1, 0, 0),
    @('Du', 2, 2, 2, 0, 1, 0, 1, 0, 0),
    @('Eder', 2, 2, 2, 1, 1, 0, 1, 0, 0),
    @('Marcelão', 2, 2, 2, 2, 1, 0, 1, 0, 0),
    @('Cabeleira', 2, 2, 2, 0, 1, 0, 1, 0, 0),
    @('Leandrão', 2, 2, 2, 2, 1, 1, 0, 0, 0),
    @('Jorge', 2, 2, 2, 0, 1, 1, 0, 0, 0),
    @('Boneco', 2, 2, 2, 3, 1, 1, 0, 0, 0),
    @('Corinthiano', 2, 2, 2, 0, 1, 1, 0, 0, 0),
    @('Athos', 2, 2, 2, 1, 1, 1, 0, 0, 0),
    @('David', 2, 2, 2, 1, 1, 0, 0, 0, 0),
    @('Philipe', 2, 2, 2, 0, 1, 0, 0, 0, 0),
    @('Juscielio', 2, 2, 2, 2, 1, 0, 0, 0, 0),
    @('Ismael', 2, 2, 2, 0, 1, 0, 0, 0, 0),
    @('Leandrinho', 2, 2, 2, 2, 1, 0, 0, 0, 0),
    @('Matheus', 4, 4, 3, 0, 1, 1, 0, 0, 10),
    @('Chelin', 3, 4, 4, 0, 1, 0, 1, 0, 11)
)

$startRow = 440
for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $startRow + $i
    $row = $data[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 3).Value = $row[1]
    $ws.Cells.Item($r, 4).Value = $row[2]
    $ws.Cells.Item($r, 5).Value = $row[3]
    $ws.Cells.Item($r, 6).Value = $row[4]
    $ws.Cells.Item($r, 7).Value = $row[5]
    $ws.Cells.Item($r, 8).Value = $row[6]
    $ws.Cells.Item($r, 9).Value = $row[7]
    $ws.Cells.Item($r, 10).Value = $row[8]
    $ws.Cells.Item($r, 11).Value = $row[9]
}

$nextRow = $startRow + $data.Count
$ws.Range("A$nextRow").Select()
